$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "All Published Values": append a new data row (row 12) captured from the
#    BOC USD rates scrape.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("All Published Values")

$newRow = @{
    A = "2026-01-02"
    B = "2026-01-02 20:28:09"
    C = "697.85"
    D = "697.85"
    E = "700.79"
    F = "700.79"
    G = "702.88"
    H = "2026/01/02 20:28:09"
    I = "2026-01-02 12:30:28"
    J = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
}

$targetRow = 12
$rowRange = $ws1.Range("A" + $targetRow + ":J" + $targetRow)

# Force the new cells to be stored as plain text (matching every other row in
# the sheet) instead of having Excel auto-infer numbers/dates from the
# look-alike strings.
$rowRange.NumberFormat = "@"

foreach ($col in @("A","B","C","D","E","F","G","H","I","J")) {
    $ws1.Range($col + $targetRow).Value = $newRow[$col]
}

# Drop back to the workbook's default (unstyled) cell style now that the
# values are locked in as text - keeps the new row's formatting identical to
# the existing data rows (no style index attached).
$rowRange.Style = "Normal"

# Re-apply the AutoFilter over the full data range, including the new row.
$ws1.AutoFilterMode = $false
[void]$ws1.Range("A1:J" + $targetRow).AutoFilter()

# The AutoFilter action above drives the AutoFilter XML, but Excel's hidden
# _FilterDatabase defined name for this sheet also needs to track the new
# range explicitly.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$" + $targetRow
    }
}

# ---------------------------------------------------------------------------
# 2. "Daily Summary": the day's publish count went from 10 to 11 now that a
#    new publish was captured.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B4").Value = 11

Write-Output "BOC USD rates updated"
